# Apply edits described by the diff:
# - Replace faculty/course related strings so the composite key no longer
#   repeats across rows (integrity constraint fix for WeeklyAttendance).
# - Update row 3..6 values accordingly.
# - Move the active cell selection from A6 to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Week 1) stays conceptually the same but the Srn text changes
# (1ru233 -> 1ru353) because the shared string itself was edited.
$ws.Range("H2").Value = "1ru353"

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "A"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "BE"
$ws.Range("F3").Value = "10cs42"
$ws.Range("G3").Value = 2014
$ws.Range("H3").Value = "1ru353"

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = "A"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "BE"
$ws.Range("F4").Value = "10cs48"
$ws.Range("G4").Value = 2014
$ws.Range("H4").Value = "1ru353"

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("C5").Value = "A"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "BE"
$ws.Range("F5").Value = "10cs49"
$ws.Range("H5").Value = "1ru353"

# Row 6
$ws.Range("A6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "BE"
$ws.Range("F6").Value = "10cs50"
$ws.Range("G6").Value = 2014
$ws.Range("H6").Value = "1ru353"

# Move the selection to A3 (was A6)
$ws.Range("A3").Select() | Out-Null
